$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Tipo" column (D), shifting it to E.
$ws.Range("D1").EntireColumn.Insert()

# New header for inserted column; it inherits the bordered/bold header
# style from the neighboring header cells automatically on column insert.
$ws.Range("D1").Value = "MAE"

# New MAE value for row 2.
$ws.Range("D2").Value = 0.2144501238200131

# Updated MSE (B2) and R2 (C2) values.
$ws.Range("B2").Value = 0.08041674500432616
$ws.Range("C2").Value = 0.9992372111350304
